# "Generate Report for Handback"
#
# The handback transform failed (file name mismatch between the handback
# package and what was handed off), so this report needs to reflect the
# failure: the Status cells flip from "Ready for handoff" to
# "Handback transform failed", and the (previously empty) Error Detail
# cells get a diagnostic message for each locale. A few report columns
# are also widened so the new, longer text is readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handback transform failed" ---
# (Overview rolls up the per-locale status in columns E/F, the per-locale
#  sheets carry it in column C.)
$wsOverview.Range("E2").Value = "Handback transform failed"
$wsOverview.Range("F2").Value = "Handback transform failed"
$wsZhCn.Range("C2").Value = "Handback transform failed"
$wsDeDe.Range("C2").Value = "Handback transform failed"

# --- Error Detail: explain why the handback transform failed ---
$wsZhCn.Range("P2").Value = "Handback file name: pqtubc4s.jnx is different with handoff file name: 489f0778-1eee-4465-a576-e978fbbd88a9.de4a393c341c76c61eac25e4facb5b38b5e31417.zh-cn."
$wsDeDe.Range("P2").Value = "Handback file name: pqtubc4s.jnx is different with handoff file name: 489f0778-1eee-4465-a576-e978fbbd88a9.de4a393c341c76c61eac25e4facb5b38b5e31417.de-de."

# --- Widen columns so the longer Status / Error Detail text is legible ---
# (ColumnWidth is internally snapped to a whole-pixel grid by Excel, same
#  as real Excel's automation model; the inputs below are chosen so the
#  stored width lands on (or as close as the grid allows to) the target.)
$wsOverview.Columns.Item(5).ColumnWidth = 23.833333333333336   # E -> ~24.74
$wsOverview.Columns.Item(6).ColumnWidth = 23.833333333333336   # F -> ~24.74

$wsZhCn.Columns.Item(3).ColumnWidth  = 23.833333333333336      # C (Status) -> ~24.74
$wsZhCn.Columns.Item(16).ColumnWidth = 39.16666666666667        # P (Error Detail) -> 40

$wsDeDe.Columns.Item(3).ColumnWidth  = 23.833333333333336      # C (Status) -> ~24.74
$wsDeDe.Columns.Item(16).ColumnWidth = 39.16666666666667        # P (Error Detail) -> 40
